# Applies the "Added more options for creation of design matrix, and senstype ref."
# commit to the design_input_example1 workbook.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("general_input")
$wsDesign  = $wb.Worksheets.Item("designinput")
$wsDefault = $wb.Worksheets.Item("defaultvalues")

# ---------------------------------------------------------------------------
# 1. designinput ("designinput" sheet) header + example-row relabeling
# ---------------------------------------------------------------------------

# Header row: senscase1 / senscase2 -> casename1 / casename2
$wsDesign.Range("E1").Value = "casename1"
$wsDesign.Range("G1").Value = "casename2"

# Row 2 (seed example row): rename sensname label, fix type, clear stray
# param_name value that leaked into the example row.
$wsDesign.Range("A2").Value = "rms_seed"
$wsDesign.Range("C2").Value = "seed"
$wsDesign.Range("D2").ClearContents()

# Row 3-8 param_name labels -> upper-cased parameter identifiers
$wsDesign.Range("D3").Value = "FAULT_POSITION"
$wsDesign.Range("D4").Value = "DC_MODEL"
$wsDesign.Range("D5").Value = "OWC1"
$wsDesign.Range("D6").Value = "OWC2"
$wsDesign.Range("D7").Value = "OWC3"
$wsDesign.Range("D8").Value = "MULTZ_ILE"

# Alignment fix: E3 / E4 / E5 become left aligned (were general/center/general)
$wsDesign.Range("E3").HorizontalAlignment = -4131
$wsDesign.Range("E4").HorizontalAlignment = -4131
$wsDesign.Range("E5").HorizontalAlignment = -4131

# Column D is widened to fit the longer upper-case parameter names
# (target character width ~15.19; engine quantizes to pixel steps so this is
# the closest achievable value)
$wsDesign.Columns.Item(4).ColumnWidth = 14.33

# ---------------------------------------------------------------------------
# 2. Comment text updates on designinput!E1 and designinput!I1
# ---------------------------------------------------------------------------

$commentE1 = $wsDesign.Range("E1").Comment
$commentE1.Text("For scenario sensitivities provide names for case1 and case2 and values. Values can be string or numbers") | Out-Null

$commentI1 = $wsDesign.Range("I1").Comment
$commentI1.Text("Distname and dist_param1, .. only for sensitivities of type " + [char]0x201C + "dist" + [char]0x201D + ". " + [char]10 + "The order of distribution parameters is predefined: " + [char]10 + "normal(mean, std dev,min, max)     " + [char]0x2013 + " where min/max is optional and will give truncated gaussian" + [char]10 + "lognormal(mean, stddev) " + [char]10 + "uniform(from,to)" + [char]10 + "loguniform(from, to)" + [char]10 + "triangular(low, mode, high)" + [char]10 + "discrete((value1, value2, value3,..,value_n) (weight1, weight2, weight3,..weight_n)). Discrete uniform if no weights are given") | Out-Null

# ---------------------------------------------------------------------------
# 3. defaultvalues sheet: keep param_name column in sync with designinput,
#    and relabel the base-case parameter placeholders.
# ---------------------------------------------------------------------------

$wsDefault.Range("A3").Value = "FAULT_POSITION"
$wsDefault.Range("A4").Value = "DC_MODEL"
$wsDefault.Range("A5").Value = "OWC1"
$wsDefault.Range("A6").Value = "OWC2"
$wsDefault.Range("A7").Value = "OWC3"
$wsDefault.Range("A8").Value = "MULTZ_ILE"
$wsDefault.Range("A9").Value = "PARAM1"
$wsDefault.Range("A10").Value = "PARAM2"
$wsDefault.Range("A11").Value = "PARAM3"
$wsDefault.Range("A12").Value = "PARAM4"

# ---------------------------------------------------------------------------
# 4. Active-sheet / selection bookkeeping: designinput becomes the active tab
#    (was defaultvalues), with its own new selection, while defaultvalues'
#    and general_input's stored selections move too.
# ---------------------------------------------------------------------------

$wsDefault.Range("B12").Select() | Out-Null

$wsDesign.Activate() | Out-Null
$wsDesign.Range("Q35").Select() | Out-Null
